# Apply the updated cryptocurrency price/volume snapshot (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.828.74"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3
$ws.Range("E3").Value = "  +1.31%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "'245.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.50%  "

# Row 6
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").Value = "'0.4790"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.86%  "

# Row 8
$ws.Range("D8").Value = "'0.2910"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.31%  "

# Row 9
$ws.Range("D9").Value = "'43.10"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.42%  "

# Row 10
$ws.Range("D10").Value = "'0.06576"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.41%  "

# Row 11
$ws.Range("D11").Value = "'21.38"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.93%  "

# Row 12
$ws.Range("D12").Value = "'0.07785"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "

# Row 13
$ws.Range("D13").Value = "1.910.75"
$ws.Range("E13").Value = "  +2.92%  "

# Row 14
$ws.Range("D14").Value = "'97.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.77%  "

# Row 15
$ws.Range("D15").Value = "'0.7425"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.68%  "

# Row 16
$ws.Range("D16").Value = "'5.193"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.16%  "

# Row 17
$ws.Range("D17").Value = "'281.85"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.79%  "

# Row 18
$ws.Range("D18").Value = "30.820.92"
$ws.Range("E18").Value = "  +2.01%  "

# Row 19
$ws.Range("D19").Value = "'13.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.23%  "

# Row 20
$ws.Range("D20").Value = "'0.000007611"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.15%  "

# Row 21
$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.01%  "

# Row 22
$ws.Range("D22").Value = "2.156.65"
$ws.Range("E22").Value = "  +2.04%  "

# Row 23
$ws.Range("D23").Value = "'5.316"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.63%  "

# Row 24
$ws.Range("E24").Value = "  +0.02%  "

# Row 25
$ws.Range("D25").Value = "'6.252"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.60%  "

# Row 26
$ws.Range("D26").Value = "'9.380"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.21%  "

# Row 27
$ws.Range("D27").Value = "'166.58"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.60%  "

# Row 28
$ws.Range("D28").Value = "'19.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.45%  "

# Row 29
$ws.Range("D29").Value = "'1.984"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.78%  "

# Row 30
$ws.Range("B30").Value = "Stellar"
$ws.Range("C30").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D30").Value = "'0.1004"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.23%  "

# Row 31
$ws.Range("B31").Value = "Toncoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D31").Value = "'1.374"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.25%  "

# Row 32
$ws.Range("D32").Value = "'1.519"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.56%  "

# Row 33
$ws.Range("D33").Value = "'4.389"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.05%  "

# Row 34
$ws.Range("D34").Value = "'4.138"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.44%  "

# Row 35
$ws.Range("D35").Value = "'0.04797"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.41%  "

# Row 36
$ws.Range("E36").Value = "  +0.52%  "

# Row 37
$ws.Range("D37").Value = "'0.7062"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.00%  "

# Row 38
$ws.Range("D38").Value = "'2.719"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.04%  "

# Row 39
$ws.Range("E39").Value = "  +0.89%  "

# Row 40
$ws.Range("D40").Value = "'2.771"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "

# Row 41
$ws.Range("D41").Value = "'6.443"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.00%  "

# Row 42
$ws.Range("D42").Value = "'70.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.26%  "

# Row 43
$ws.Range("D43").Value = "'1.936"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.01%  "

# Row 44
$ws.Range("D44").Value = "'0.4223"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.75%  "

# Row 45
$ws.Range("D45").Value = "'0.8490"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "

# Row 46
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'102.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48
$ws.Range("D48").Value = "'9.513"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.10%  "

# Row 49
$ws.Range("D49").Value = "'7.171"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.39%  "

# Row 50
$ws.Range("D50").Value = "'946.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.17%  "

# Row 51
$ws.Range("D51").Value = "'35.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.73%  "
